# Updated cryptos list (Price / Volume(1h) columns) on the single worksheet.
# Numeric-looking text values in column D are prefixed with a leading
# apostrophe so Excel keeps them as text (matching the original cell
# content type) instead of auto-converting them to numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "44.187.79"
$ws.Range("E2").Value = "  +5.28%  "
$ws.Range("D3").Value = "2.266.49"
$ws.Range("E3").Value = "  +2.78%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").Value = "'230.62"
$ws.Range("E5").Value = "  +0.04%  "
$ws.Range("D6").Value = "'0.634"
$ws.Range("E6").Value = "  +2.86%  "
$ws.Range("D7").Value = "'63.71"
$ws.Range("E7").Value = "  +5.23%  "
$ws.Range("E8").Value = "  +0.10%  "
$ws.Range("D9").Value = "'0.448"
$ws.Range("E9").Value = "  +11.82%  "
$ws.Range("D10").Value = "'0.105"
$ws.Range("E10").Value = "  +16.57%  "
$ws.Range("D11").Value = "'56.88"
$ws.Range("E11").Value = "  -0.50%  "
$ws.Range("D12").Value = "'26.40"
$ws.Range("E12").Value = "  +20.34%  "
$ws.Range("E13").Value = "  +2.95%  "
$ws.Range("D14").Value = "2.601.04"
$ws.Range("E14").Value = "  +2.86%  "
$ws.Range("D15").Value = "'15.76"
$ws.Range("E15").Value = "  +2.56%  "
$ws.Range("E16").Value = "  +9.46%  "
$ws.Range("D17").Value = "'0.836"
$ws.Range("D18").Value = "2.255.81"
$ws.Range("E18").Value = "  +2.33%  "
$ws.Range("D19").Value = "43.994.88"
$ws.Range("E19").Value = "  +5.07%  "
$ws.Range("D20").Value = "'0.0000102"
$ws.Range("E20").Value = "  +9.05%  "
$ws.Range("D21").Value = "'73.54"
$ws.Range("E21").Value = "  +2.49%  "
$ws.Range("D22").Value = "'6.04"
$ws.Range("E22").Value = "  -2.10%  "
$ws.Range("D23").Value = "'252.08"
$ws.Range("E24").Value = "  +0.12%  "
$ws.Range("E25").Value = "  +1.31%  "
$ws.Range("E26").Value = "  -1.60%  "
$ws.Range("E27").Value = "  +27.58%  "
$ws.Range("D28").Value = "'10.07"
$ws.Range("E28").Value = "  +5.08%  "
$ws.Range("D29").Value = "'171.78"
$ws.Range("E29").Value = "  +1.79%  "
$ws.Range("D30").Value = "'20.81"
$ws.Range("E30").Value = "  +3.03%  "
$ws.Range("E31").Value = "  -1.25%  "
$ws.Range("E32").Value = "  -2.97%  "
$ws.Range("E33").Value = "  +3.55%  "
$ws.Range("D34").Value = "'0.0683"
$ws.Range("E34").Value = "  +6.11%  "
$ws.Range("D35").Value = "'4.76"
$ws.Range("E35").Value = "  +3.85%  "
$ws.Range("D36").Value = "'4.88"
$ws.Range("E36").Value = "  -0.83%  "
$ws.Range("D37").Value = "'3.85"
$ws.Range("E37").Value = "  +9.44%  "
$ws.Range("D38").Value = "'6.70"
$ws.Range("E38").Value = "  +6.84%  "
$ws.Range("E39").Value = "  -0.19%  "
$ws.Range("D40").Value = "'0.0259"
$ws.Range("E40").Value = "  +4.84%  "
$ws.Range("E41").Value = "  +0.18%  "
$ws.Range("D42").Value = "'17.55"
$ws.Range("E42").Value = "  +9.48%  "
$ws.Range("D43").Value = "'8.32"
$ws.Range("E43").Value = "  -2.37%  "
$ws.Range("D44").Value = "'0.0965"
$ws.Range("E44").Value = "  +1.56%  "
$ws.Range("D45").Value = "'97.83"
$ws.Range("E45").Value = "  +1.48%  "
$ws.Range("D46").Value = "'0.000212"
$ws.Range("E46").Value = "  -5.33%  "
$ws.Range("E47").Value = "  -0.43%  "
$ws.Range("E48").Value = "  +2.23%  "
$ws.Range("D49").Value = "1.446.05"
$ws.Range("E49").Value = "  -0.59%  "
$ws.Range("D50").Value = "'9.97"
$ws.Range("E50").Value = "  +19.02%  "
$ws.Range("E51").Value = "  +4.96%  "
